# edit.ps1 — applies "getting rid of duplicate actions" changes to Methods.docx
#
# Summary of changes (per the commit's xml diff):
#   1. Paragraph "Intersected with country polygon" — drop the stray
#      leading space-only run.
#   2. Paragraph "Point files combined and buffered..." — append a new
#      " (convert from km to m)" run and move the lone "_GoBack" bookmark
#      here (it previously sat alone in its own trailing paragraph).
#   3. The UNESCO-exclusion paragraph — expand the single
#      `INT_CRIT == "Not Applicable"` run into the `INT_CRIT %in%
#      c("Not Applicable", "Not Reported")` run sequence (with the
#      gramStart/gramEnd proofErr markers around "c("), and collapse the
#      three STATUS runs (with their gramStart/gramEnd proofErr markers)
#      into a single run.
#   4. The now-empty paragraph that used to hold the relocated bookmark
#      becomes a bare empty paragraph.
#
# NOTE on structure: every call that ends in `.InsertXML(...)` builds its
# pkg:package string and performs the call in the *same* function
# invocation as the one that computed the target Range — round-tripping a
# Range object through a second layer of helper-function indirection
# together with deferred string concatenation confused this host's COM
# shim (the mutation silently no-op'd), so the helpers below are kept
# intentionally flat.

function Get-ParagraphByText($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.Contains($needle)) {
            return $para
        }
    }
    return $null
}

function Set-ParagraphBodyXml($doc, $para, [string]$innerParagraphXml) {
    # Replace everything in the paragraph *except* its trailing paragraph
    # mark, so pPr / the mark's own rsid attributes survive untouched.
    $rng = $para.Range
    $body = $doc.Range($rng.Start, $rng.End - 1)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $body.InsertXML($pkg)
}

function Reset-EmptyParagraph($doc, $para) {
    # Replace the paragraph *including* its paragraph mark with a bare
    # <w:p/>, so any leftover rsid/rsidR/rsidRDefault attributes on the
    # mark itself are dropped too (used once the bookmark it used to
    # carry has been relocated elsewhere).
    $rng = $para.Range
    $full = $doc.Range($rng.Start, $rng.End)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($pkg)
}

$d = $word.ActiveDocument

# --- Change 1: remove the duplicate leading-space run -----------------
$p8 = Get-ParagraphByText $d "Intersected with country polygon"
Set-ParagraphBodyXml $d $p8 '<w:p><w:r w:rsidR="001F128A"><w:t>Intersected with country polygon</w:t></w:r></w:p>'

# --- Change 2: add the " (convert from km to m)" run + relocate the ---
# --- "_GoBack" bookmark here --------------------------------------------
$p19 = Get-ParagraphByText $d "Point files combined and buffered based on their reported area"
Set-ParagraphBodyXml $d $p19 '<w:p><w:r><w:t>Point files combined and buffered based on their reported area</w:t></w:r><w:r><w:t xml:space="preserve"> (convert from km to m)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# --- Change 3: expand/collapse the INT_CRIT / STATUS runs --------------
$p21 = Get-ParagraphByText $d "Data filtered to exclude UNESCO sites"
Set-ParagraphBodyXml $d $p21 '<w:p><w:r><w:t>Data filtered to exclude UNESCO sites</w:t></w:r><w:r w:rsidR="00D94F73"><w:t xml:space="preserve"> (</w:t></w:r><w:r w:rsidR="00D94F73" w:rsidRPr="00D94F73"><w:t xml:space="preserve">INT_CRIT </w:t></w:r><w:r><w:t>%in%</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>c(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>"Not Applicable"</w:t></w:r><w:r><w:t>, “Not Reported”</w:t></w:r><w:r w:rsidR="00D94F73"><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve">, sites with a status of “Proposed” or “Not Reported” </w:t></w:r><w:r w:rsidR="00D94F73"><w:t>(STATUS %!in% c(“Proposed”, “Not Reported”))</w:t></w:r></w:p>'

# --- Change 4: the bookmark used to live alone in its own paragraph; ---
# --- now that it has moved, delete it and leave a bare empty paragraph -
$goBack = $d.Bookmarks.Item("_GoBack")
if ($goBack -ne $null) {
    $goBack.Delete()
}
# the bookmark's old paragraph is now contentless; it is the paragraph
# immediately following the UNESCO paragraph we just edited.
$unescoPara = Get-ParagraphByText $d "Data filtered to exclude UNESCO sites"
$nextPara = $unescoPara.Next()
if ($nextPara.Range.Text -eq "") {
    Reset-EmptyParagraph $d $nextPara
}

Write-Output "done"
